$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6717289090156555
$ws.Range("B1").Value = 1.01912784576416
$ws.Range("C1").Value = 2.228631019592285
$ws.Range("D1").Value = 3.95419979095459
$ws.Range("E1").Value = 1.608527421951294
